$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 2 <-> Row 4 swap of identifying data (Id, coordinates, times, biotope description, substrate info)

# Column A - Id
$ws.Range("A2").Value = 111697236
$ws.Range("A4").Value = 111697304

# Column Q - Ost
$ws.Range("Q2").Value = 373121.3523494597
$ws.Range("Q4").Value = 373090.8741807578

# Column R - Nord
$ws.Range("R2").Value = 6865443.651501717
$ws.Range("R4").Value = 6865424.499624529

# Column Z - Starttid
$ws.Range("Z2").Value = "00:00"
$ws.Range("Z4").Value = "19:00"

# Column AB - Sluttid
$ws.Range("AB2").Value = "00:00"
$ws.Range("AB4").Value = "19:00"

# Column AI - Biotop-beskrivning
$ws.Range("AI2").Value = "Tallskog. Kontinuitetsskog"
$ws.Range("AI4").Value = "Luckig tallskog. K-skog"

# Column AJ - Substratnamn (moves from row4 to row2)
$ws.Range("AJ2").Value = "tall"
$ws.Range("AJ4").Value = ""

# Column AK - Vetenskapligt Substratnamn (moves from row4 to row2)
$ws.Range("AK2").Value = "Pinus sylvestris"
$ws.Range("AK4").Value = ""

# Column AO - Substrat-beskrivning (moves from row4 to row2)
$ws.Range("AO2").Value = "Pinus sylvestris"
$ws.Range("AO4").Value = ""
